$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '80.832.51'
$ws.Range("E2").Value = '  +5.88%  '

# Row 3
$ws.Range("D3").Value = '3.240.11'
$ws.Range("E3").Value = '  +6.40%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '214.05'
$ws.Range("E5").Value = '  +7.70%  '

# Row 6
$ws.Range("D6").Value = '637.68'
$ws.Range("E6").Value = '  +3.06%  '

# Row 7
$ws.Range("D7").Value = '0.285'
$ws.Range("E7").Value = '  +36.62%  '

# Row 8
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("D9").Value = '0.608'
$ws.Range("E9").Value = '  +10.47%  '

# Row 10
$ws.Range("D10").Value = '3.234.11'
$ws.Range("E10").Value = '  +6.27%  '

# Row 11
$ws.Range("D11").Value = '0.614'
$ws.Range("E11").Value = '  +39.60%  '

# Row 12
$ws.Range("D12").Value = '0.0000276'
$ws.Range("E12").Value = '  +43.47%  '

# Row 13
$ws.Range("E13").Value = '  +3.53%  '

# Row 14
$ws.Range("D14").Value = '5.49'
$ws.Range("E14").Value = '  +5.62%  '

# Row 15
$ws.Range("D15").Value = '3.828.98'
$ws.Range("E15").Value = '  +6.27%  '

# Row 16
$ws.Range("D16").Value = '32.94'
$ws.Range("E16").Value = '  +13.95%  '

# Row 17
$ws.Range("D17").Value = '80.484.85'
$ws.Range("E17").Value = '  +5.42%  '

# Row 18
$ws.Range("D18").Value = '3.229.92'
$ws.Range("E18").Value = '  +6.25%  '

# Row 19
$ws.Range("D19").Value = '14.72'
$ws.Range("E19").Value = '  +8.71%  '

# Row 20
$ws.Range("E20").Value = '  +26.17%  '

# Row 21
$ws.Range("D21").Value = '9.45'
$ws.Range("E21").Value = '  +6.22%  '

# Row 22
$ws.Range("D22").Value = '450.04'
$ws.Range("E22").Value = '  +18.07%  '

# Row 23
$ws.Range("D23").Value = '5.42'
$ws.Range("E23").Value = '  +23.00%  '

# Row 24
$ws.Range("E24").Value = '  +12.57%  '

# Row 25
$ws.Range("D25").Value = '3.408.82'
$ws.Range("E25").Value = '  +7.26%  '

# Row 26
$ws.Range("E26").Value = '  +7.80%  '

# Row 27
$ws.Range("D27").Value = '11.09'
$ws.Range("E27").Value = '  +12.89%  '

# Row 28
$ws.Range("E28").Value = '  +19.02%  '

# Row 29
$ws.Range("E29").Value = '  +0.06%  '

# Row 30
$ws.Range("D30").Value = '9.41'
$ws.Range("E30").Value = '  +13.33%  '

# Row 31
$ws.Range("E31").Value = '  +0.02%  '

# Row 32
$ws.Range("D32").Value = '569.67'
$ws.Range("E32").Value = '  +15.39%  '

# Row 33
$ws.Range("D33").Value = '1.54'
$ws.Range("E33").Value = '  +10.00%  '

# Row 34
$ws.Range("E34").Value = '  +30.93%  '

# Row 35
$ws.Range("E35").Value = '  +6.93%  '

# Row 36
$ws.Range("D36").Value = '23.89'
$ws.Range("E36").Value = '  +15.74%  '

# Row 37
$ws.Range("E37").Value = '  +21.18%  '

# Row 38
$ws.Range("E38").Value = '  +11.19%  '

# Row 39
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.03%  '

# Row 40
$ws.Range("D40").Value = '5.88'
$ws.Range("E40").Value = '  +14.68%  '

# Row 41
$ws.Range("D41").Value = '164.67'
$ws.Range("E41").Value = '  +1.49%  '

# Row 42
$ws.Range("D42").Value = '20.35'
$ws.Range("E42").Value = '  +1.44%  '

# Row 43
$ws.Range("D43").Value = '193.76'
$ws.Range("E43").Value = '  +0.81%  '

# Row 45
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.79'
$ws.Range("E45").Value = '  +14.67%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '1.86'
$ws.Range("E46").Value = '  +13.36%  '

# Row 47
$ws.Range("E47").Value = '  +9.57%  '

# Row 48
$ws.Range("D48").Value = '0.807'
$ws.Range("E48").Value = '  +3.13%  '

# Row 49
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '43.72'
$ws.Range("E49").Value = '  +6.19%  '

# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '4.38'
$ws.Range("E50").Value = '  +13.45%  '

# Row 51
$ws.Range("D51").Value = '0.652'
$ws.Range("E51").Value = '  +9.79%  '
